$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: shift Date_Creation/Nb_Produits/Valeur_Stock_Total columns
$ws.Range("H1").Value = "Nb_Produits"
$ws.Range("I1").Value = "Valeur_Stock_Total"
$ws.Range("J1").Value = "Date_Creation"

# Update row 2 values
$ws.Range("B2").Value = "Fournisseur C"
$ws.Range("F2").Value = "À définir"
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0

# Force J2 to stay a plain text value ("2025-06-04") instead of being
# auto-converted into a date serial number, then drop the temporary
# text number-format so the cell keeps the default (no) style.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2025-06-04"
$ws.Range("J2").ClearFormats()

# Delete rows 3-5 (old data rows no longer present)
$ws.Range("A3:J5").Delete()
